$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 update: report period moved from 2017-12-31 to 2018-12-31,
# along with refreshed cash-flow figures for the new period.
$ws.Range("N2").Value = "2018-12-31 00:00:00"

$ws.Range("O2").Value = 7080951.17
$ws.Range("P2").Value = 175.2828316043
$ws.Range("Q2").Value = 315550652.1
$ws.Range("R2").Value = 7811.1839054902
$ws.Range("S2").Value = 67325897.91
$ws.Range("T2").Value = 1666.5944648741
$ws.Range("U2").Value = -32958304.71
$ws.Range("V2").Value = -815.854372633
$ws.Range("Y2").Value = 33081926.21
$ws.Range("Z2").Value = 818.9145161148
$ws.Range("AA2").Value = 27274747.7
$ws.Range("AB2").Value = 675.1628267687
$ws.Range("AC2").Value = 4039728.88
$ws.Range("AD2").Value = 109.3773229227
